$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FSPbPPT")

$newItems = @(
    "hard coal w CCS",
    "natural gas combined cycle w CCS",
    "biomass w CCS",
    "lignite w CCS",
    "small modular reactor",
    "hydrogen"
)

$row = 19
foreach ($item in $newItems) {
    $ws.Cells.Item($row, 1).Value = $item
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value = 0
    $cell.NumberFormat = "0"
    $row = $row + 1
}

$ws.Range("A25").Select()
